# Auto-generated script applying the market-data refresh diff
# (currentAveragePrice / LevePrice / LeveProfit columns H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 602.25
$ws.Range("I2").Value = 480
$ws.Range("J2").Value = 724.5
$ws.Range("K2").Value = 480
$ws.Range("L2").Value = 724.5
$ws.Range("M2").Value = -367
$ws.Range("N2").Value = -950.5
# Row 6
$ws.Range("H6").Value = 125083.625
$ws.Range("I6").Value = 125083.625
$ws.Range("K6").Value = 375250.875
$ws.Range("M6").Value = -375138.875
# Row 43
$ws.Range("H43").Value = 5283.5
$ws.Range("I43").Value = 4625.75
$ws.Range("K43").Value = 4625.75
$ws.Range("M43").Value = -4556.75
# Row 46
$ws.Range("H46").Value = 1898.3334
$ws.Range("I46").Value = 1898.3334
$ws.Range("K46").Value = 5695.0002
$ws.Range("M46").Value = -5576.0002
# Row 48
$ws.Range("H48").Value = 7833.3335
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 30000
$ws.Range("M48").Value = -29708
# Row 56
$ws.Range("H56").Value = 7833.3335
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 30000
$ws.Range("M56").Value = -29466
# Row 60
$ws.Range("H60").Value = 1898.3334
$ws.Range("I60").Value = 1898.3334
$ws.Range("K60").Value = 5695.0002
$ws.Range("M60").Value = -5211.0002
# Row 74
$ws.Range("H74").Value = 8591.857
$ws.Range("I74").Value = 8591.857
$ws.Range("K74").Value = 8591.857
$ws.Range("M74").Value = -7655.857
# Row 77
$ws.Range("H77").Value = 8591.857
$ws.Range("I77").Value = 8591.857
$ws.Range("K77").Value = 42959.285
$ws.Range("M77").Value = -38279.285
# Row 86
$ws.Range("H86").Value = 3918
$ws.Range("I86").Value = 3928.7273
$ws.Range("K86").Value = 3928.7273
$ws.Range("M86").Value = -2805.7273
# Row 89
$ws.Range("H89").Value = 3918
$ws.Range("I89").Value = 3928.7273
$ws.Range("K89").Value = 19643.6365
$ws.Range("M89").Value = -14027.6365
# Row 92
$ws.Range("H92").Value = 511.33334
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 100
$ws.Range("H100").Value = 1847.6666
$ws.Range("I100").Value = 2163.9167
$ws.Range("K100").Value = 2163.9167
$ws.Range("M100").Value = -1622.9167

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 888
$ws.Range("I4").Value = 888
$ws.Range("K4").Value = 888
$ws.Range("M4").Value = -772
# Row 32
$ws.Range("H32").Value = 15552.318
$ws.Range("I32").Value = 16429.264
$ws.Range("K32").Value = 16429.264
$ws.Range("M32").Value = -16142.264
# Row 45
$ws.Range("H45").Value = 1168.3334
$ws.Range("I45").Value = 1168.3334
$ws.Range("K45").Value = 1168.3334
$ws.Range("M45").Value = -791.3334
# Row 58
$ws.Range("H58").Value = 55555
$ws.Range("J58").Value = 55555
$ws.Range("L58").Value = 55555
$ws.Range("N58").Value = -56415
# Row 61
$ws.Range("H61").Value = 1699.5
$ws.Range("I61").Value = 1699.5
$ws.Range("K61").Value = 1699.5
$ws.Range("M61").Value = -1487.5
# Row 74
$ws.Range("H74").Value = 3192.4211
$ws.Range("I74").Value = 1822.8
$ws.Range("J74").Value = 4714.222
$ws.Range("K74").Value = 1822.8
$ws.Range("L74").Value = 4714.222
$ws.Range("M74").Value = -948.8
$ws.Range("N74").Value = -6462.222
# Row 76
$ws.Range("H76").Value = 79949.5
$ws.Range("J76").Value = 79949.5
$ws.Range("L76").Value = 79949.5
$ws.Range("N76").Value = -80625.5
# Row 77
$ws.Range("H77").Value = 3192.4211
$ws.Range("I77").Value = 1822.8
$ws.Range("J77").Value = 4714.222
$ws.Range("K77").Value = 9114
$ws.Range("L77").Value = 23571.11
$ws.Range("M77").Value = -4746
$ws.Range("N77").Value = -32307.11
# Row 79
$ws.Range("H79").Value = 79949.5
$ws.Range("J79").Value = 79949.5
$ws.Range("L79").Value = 79949.5
$ws.Range("N79").Value = -82289.5
# Row 132
$ws.Range("H132").Value = 3693.0571
$ws.Range("I132").Value = 2641.9333
$ws.Range("K132").Value = 7925.7999
$ws.Range("M132").Value = -5395.7999
# Row 136
$ws.Range("H136").Value = 1699.5
$ws.Range("I136").Value = 1699.5
$ws.Range("K136").Value = 5098.5
$ws.Range("M136").Value = -2548.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 9163.429
$ws.Range("I7").Value = 1318
$ws.Range("J7").Value = 28777
$ws.Range("K7").Value = 1318
$ws.Range("L7").Value = 28777
$ws.Range("M7").Value = -1205
$ws.Range("N7").Value = -29003

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 2349.5715
$ws.Range("I17").Value = 1741.1666
$ws.Range("K17").Value = 1741.1666
$ws.Range("M17").Value = -1567.1666
# Row 25
$ws.Range("H25").Value = 6391.1665
$ws.Range("I25").Value = 4669.4
$ws.Range("K25").Value = 4669.4
$ws.Range("M25").Value = -4495.4
# Row 50
$ws.Range("H50").Value = 26381.5
$ws.Range("I50").Value = 15210.4
$ws.Range("K50").Value = 15210.4
$ws.Range("M50").Value = -14585.4
# Row 105
$ws.Range("H105").Value = 4000
$ws.Range("J105").Value = 4000
$ws.Range("L105").Value = 4000
$ws.Range("N105").Value = -7494

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 43859.72
$ws.Range("J34").Value = 47628
$ws.Range("L34").Value = 142884
$ws.Range("N34").Value = -143052
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 113
$ws.Range("H113").Value = 2068.1904
$ws.Range("J113").Value = 1913.5454
$ws.Range("L113").Value = 5740.6362
$ws.Range("N113").Value = -10080.6362
# Row 121
$ws.Range("H121").Value = 3288.3635
$ws.Range("I121").Value = 550
$ws.Range("J121").Value = 4315.25
$ws.Range("K121").Value = 1650
$ws.Range("L121").Value = 12945.75
$ws.Range("M121").Value = -340
$ws.Range("N121").Value = -15565.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1257.5
$ws.Range("I97").Value = 710
$ws.Range("J97").Value = 2900
$ws.Range("K97").Value = 710
$ws.Range("L97").Value = 2900
$ws.Range("M97").Value = -214
$ws.Range("N97").Value = -3892
# Row 102
$ws.Range("H102").Value = 1618.16
$ws.Range("I102").Value = 1268.9166
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 1268.9166
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = 353.0834
$ws.Range("N102").Value = -13244
# Row 126
$ws.Range("H126").Value = 3064.6
$ws.Range("I126").Value = 2791.0715
$ws.Range("K126").Value = 8373.2145
$ws.Range("M126").Value = -5903.2145
# Row 127
$ws.Range("H127").Value = 75000
$ws.Range("J127").Value = 75000
$ws.Range("L127").Value = 75000
$ws.Range("N127").Value = -84920
# Row 132
$ws.Range("H132").Value = 9261664
$ws.Range("I132").Value = 2115.4482
$ws.Range("K132").Value = 6346.344599999999
$ws.Range("M132").Value = -3816.344599999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# Row 40
$ws.Range("H40").Value = 200006580
$ws.Range("I40").Value = 500002500
$ws.Range("J40").Value = 9299
$ws.Range("K40").Value = 500002500
$ws.Range("L40").Value = 9299
$ws.Range("M40").Value = -500002364
$ws.Range("N40").Value = -9571

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 16249.75
$ws.Range("I26").Value = 13333
$ws.Range("K26").Value = 13333
$ws.Range("M26").Value = -13040
# Row 62
$ws.Range("H62").Value = 5199.75
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 5899.5
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 5899.5
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -7147.5
# Row 65
$ws.Range("H65").Value = 5199.75
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 5899.5
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 29497.5
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -35737.5
# Row 121
$ws.Range("H121").Value = 100000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 100000
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 100000
$ws.Range("N121").Value = -103494
$ws.Range("M121").ClearContents()
# Row 122
$ws.Range("H122").Value = 2143
$ws.Range("I122").Value = 2143
$ws.Range("K122").Value = 6429
$ws.Range("M122").Value = -3979
# Row 125
$ws.Range("H125").Value = 29998
$ws.Range("J125").Value = 29998
$ws.Range("L125").Value = 29998
$ws.Range("N125").Value = -39838
# Row 126
$ws.Range("H126").Value = 800.25
$ws.Range("I126").Value = 800.25
$ws.Range("K126").Value = 2400.75
$ws.Range("M126").Value = 69.25

Write-Host "Applied Phantom_Profits market-data refresh"
